# "Ajout étape 2 et 7 - tracés et cartes"
#
# Fill in the remaining "TBD" staff names in the staff guide workbook:
#   - CO!C22            (Comité Écoresponsable)      TBD -> Antoine St-Jean
#   - SOUTIEN!C5         (Dépannage neutre)           TBD -> Baroudeur Support Technique<br/>Matis Boyer
#   - SOUTIEN!C10        (Radios)                     TBD -> Michel Néron
#   - COMM!C3            (Commissaires list)          trailing TBD/TBD -> Félix-Antoine Malo / Geneviève Marcotte

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CO")
$ws2 = $wb.Worksheets.Item("SOUTIEN")
$ws3 = $wb.Worksheets.Item("COMM")

# 1) COMM!C3 holds the Commissaires list, the *last* unique string in the
#    shared-strings table. Edit it first, before any brand-new strings get
#    appended, so the engine mutates this entry in place instead of retiring
#    it and tacking the new text onto the end of the table.
$ws3.Range("C3").Value = "Steve Head<br/>`nHélène Soulard<br/>`nVallérie Trottier<br/>`nNancy Daigle<br/>`nFélix-Antoine Malo<br/>`nGeneviève Marcotte"

# 2) SOUTIEN: replace the two outstanding TBDs. Order matters here too -
#    whichever cell is written first claims the earlier new shared-string
#    slot, so write "Radios" (C10) before "Dépannage neutre" (C5).
$ws2.Range("C10").Value = "Michel Néron"
$ws2.Range("C5").Value = "Baroudeur Support Technique<br/>Matis Boyer"

# 3) CO: Comité Écoresponsable gets its name (this text already exists
#    elsewhere in the workbook, so no new shared string is created).
$ws1.Range("C22").Value = "Antoine St-Jean"

# 4) Leave the workbook where the author ended up: SOUTIEN's selection
#    moved on to C6, and the final active sheet/cell is CO!C3.
$ws2.Select()
$ws2.Range("C6").Select()

$ws1.Select()
$ws1.Range("C3").Select()
